$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: C, D, E, F, G, H, I, K, N, O (as numeric indices)
$cols = @(3,4,5,6,7,8,9,11,14,15)

$data = @(
    @(0.01867454579619476,0.01614318093050571,0.4248923627222041,0.4087128934447577,0.2599563856376221,0.4151933184147083,0.2983423642013037,1.02386214922933,0.8913632644226368,1.269569322857564),
    @(0.01628825686935187,0.01419227114529775,0.3706619140378393,0.4031517034930729,0.2555352550209236,0.4170207067384695,0.2948900495461544,0.8933833986952209,0.883633951762377,1.263743596995951),
    @(0.01481989529445116,0.01298827817062431,0.3374475790228502,0.4000820561260952,0.2530913040299154,0.418397931638367,0.2930257512682708,0.8129774182268363,0.879216026500103,1.261264851146535),
    @(0.01422076225017577,0.01249613865913091,0.3239318125341839,0.3989176585241054,0.2521630911586499,0.4190232619484391,0.2923300822973651,0.7801401465794697,0.8774985576277317,1.260530037131318),
    @(0.01412123152073264,0.01241432963453093,0.3216886509361387,0.3987295293492394,0.2520130421681941,0.4191309669431647,0.2922184299393855,0.7746833059486846,0.8772183908853037,1.260424623385106),
    @(0.01481181822566668,0.01298164703935356,0.33726522447202,0.4000660026899823,0.2530785120921379,0.4184061056211519,0.293016110151413,0.8125348486962309,0.879192528013661,1.261253827693778),
    @(0.01785243077317489,0.01547180094174649,0.4061754055134656,0.4067236598048964,0.2583756028297159,0.4157703905665073,0.2970988774757757,0.9789348063767136,0.8886303847752259,1.267332133498826),
    @(0.02378871772459945,0.02030493194715177,0.5420490652899019,0.4225280545349648,0.2709267809955804,0.4126303027463791,0.3071411571075018,1.302854067094017,0.9097234941101817,1.288008208905325),
    @(0.02813286610366106,0.02382371202507017,0.6424461764504201,0.4358334763099379,0.2814908754358925,0.4115654997890061,0.3157744682862074,1.539298622625722,0.9267779612803935,1.308599588287166),
    @(0.03010516242193262,0.02541720200218123,0.688269653457624,0.4422583125526955,0.2865935063286713,0.4113520737260643,0.3199776526202669,1.646513399764501,0.9348709719671717,1.319153535902927),
    @(0.03085143343459151,0.02601954281016106,0.7056456657216046,0.4447450282200123,0.2885688342057193,0.4113103122304977,0.3216091793221096,1.687061401549784,0.9379833884116096,1.323321736266479),
    @(0.03069073766286579,0.02588986659105075,0.7019023504318227,0.4442070739960968,0.2881414904914692,0.4113175676704373,0.3212560242668872,1.678331008840701,0.9373109553000205,1.322416389887223),
    @(0.03016657071927398,0.02546677880750536,0.6896987042562586,0.4424618170515799,0.2867551524007581,0.4113478546506997,0.3201110790333246,1.64985036166604,0.9351260771394152,1.3194930104994),
    @(0.02984542477278751,0.02520748355932767,0.6822267515494218,0.441399806801364,0.2859116004987783,0.4113714956712045,0.3194149656254126,1.632398315196212,0.9337939859087783,1.317724739289673),
    @(0.0280038906861364,0.02371942490525214,0.6394547190730577,0.4354211084283861,0.2811634138078887,0.4115849081546941,0.3155053486337565,1.532284745637014,0.9262557691776436,1.307933829179717),
    @(0.02687314939397822,0.02280467054220026,0.613255718895374,0.4318488660811823,0.2783268785977953,0.4117852960428507,0.3131777227167092,1.470778307306318,0.9217168057991501,1.30223209849575),
    @(0.02622241405424575,0.02227785081413458,0.5982009791838294,0.4298292183375523,0.2767233089191876,0.4119260520110828,0.3118648856791921,1.435369090553309,0.9191376539012737,1.299064264528681),
    @(0.0260020250583608,0.0220993636432425,0.5931060942490234,0.4291514057396242,0.2761851536228903,0.4119780860159636,0.3114248326057805,1.423374653709118,0.9182698261406728,1.298010836488828),
    @(0.02699355655758495,0.02290211809141596,0.616043161580393,0.432225511944452,0.2786259392582195,0.4117613247583165,0.3134228150627152,1.477329136389301,0.9221967247386544,1.302827494887993),
    @(0.03032054769494152,0.02559107956121665,0.6932825549099988,0.4429729799809792,0.28716118159997,0.411337897840454,0.3204462934508356,1.65821724508595,0.9357665353852269,1.320347011516958),
    @(0.03249144460961872,0.02734216167407055,0.7439016723500202,0.4503106012597442,0.2929906849181236,0.411288870027775,0.3252690424205724,1.77613446857805,0.9449134534473842,1.332797990394198),
    @(0.03133312710532721,0.02640816701085669,0.7168720222468039,0.4463656017903617,0.2898562602169932,0.411294170501165,0.3226737089515979,1.713228362169446,0.9400062342076438,1.326060754868479),
    @(0.02693912252409802,0.02285806490255027,0.6147829358119452,0.4320551242044175,0.278490649258174,0.4117720825920514,0.3133119298757805,1.474367656982395,0.9219796587882314,1.302557973081662),
    @(0.02218573747741459,0.01900296729994722,0.5052000581514307,0.4179563898147549,0.2672972704664431,0.413262036916862,0.3042050346599723,1.215489722159589,0.9037424316440337,1.281470882896969)
)

$startRow = 2
$row = $startRow
foreach ($rowVals in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $rowVals[$i]
    }
    $row++
}
